$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (avoid Excel auto-numeric coercion)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.003.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.607.49'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.598'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.82'
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.143'
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.064.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '58.965.29'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.91'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.617.42'
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '338.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.68'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.997'
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0760'
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '154.39'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.03'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.96'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.905'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.888'
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.11'
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '283.09'
$ws.Range("D41").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.64'
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.957.79'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '118.97'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.50'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.98'
$ws.Range("D51").Style = "Normal"

# Coin name / link updates (row reorder: Maker, Aave, RenderToken)
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

# Volume(1h) percentage text updates
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("E3").Value = '  -1.22%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  +3.54%  '
$ws.Range("E6").Value = '  -1.09%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  +4.45%  '
$ws.Range("E9").Value = '  -1.61%  '
$ws.Range("E10").Value = '  -0.73%  '
$ws.Range("E11").Value = '  +5.08%  '
$ws.Range("E12").Value = '  -0.82%  '
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("E15").Value = '  -1.94%  '
$ws.Range("E16").Value = '  -1.32%  '
$ws.Range("E17").Value = '  -1.71%  '
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("E20").Value = '  -1.81%  '
$ws.Range("E21").Value = '  -0.48%  '
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("E24").Value = '  +2.90%  '
$ws.Range("E25").Value = '  -0.18%  '
$ws.Range("E26").Value = '  -2.01%  '
$ws.Range("E27").Value = '  -1.18%  '
$ws.Range("E28").Value = '  +1.58%  '
$ws.Range("E30").Value = '  +1.47%  '
$ws.Range("E31").Value = '  +1.43%  '
$ws.Range("E32").Value = '  +2.32%  '
$ws.Range("E33").Value = '  +1.06%  '
$ws.Range("E34").Value = '  -0.83%  '
$ws.Range("E35").Value = '  +7.88%  '
$ws.Range("E36").Value = '  +6.41%  '
$ws.Range("E37").Value = '  -0.33%  '
$ws.Range("E38").Value = '  -0.79%  '
$ws.Range("E39").Value = '  +1.06%  '
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("E41").Value = '  -0.95%  '
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("E45").Value = '  +1.19%  '
$ws.Range("E46").Value = '  -0.94%  '
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("E49").Value = '  +7.07%  '
$ws.Range("E50").Value = '  -1.34%  '
$ws.Range("E51").Value = '  -2.16%  '
